{"js": "// Replace the 100 arithmetic-answer cells (20 rows x 5 cols) in the single\n// table of the document with new equations, per the diff. Replacement is\n// done positionally (row-major order matches the document / diff order)\n// because a few \"old\" values repeat and a naive global find/replace would\n// mis-map those duplicates.\nconst REPLACEMENTS = [\n  [\"74-9=65\",\"56+7=63\"],\n  [\"40+30=70\",\"97-19=78\"],\n  [\"2+72=74\",\"92-20=72\"],\n  [\"83-11=72\",\"18+5=23\"],\n  [\"6+80=86\",\"24+74=98\"],\n  [\"46+18=64\",\"99-34=65\"],\n  [\"20+53=73\",\"89-63=26\"],\n  [\"89-78=11\",\"8+8=16\"],\n  [\"34-8=26\",\"9+3=12\"],\n  [\"10+40=50\",\"74-8=66\"],\n  [\"33+58=91\",\"81-23=58\"],\n  [\"64+17=81\",\"22+0=22\"],\n  [\"15-14=1\",\"43-33=10\"],\n  [\"86-4=82\",\"96-69=27\"],\n  [\"46+46=92\",\"61+23=84\"],\n  [\"53-4=49\",\"32+53=85\"],\n  [\"55-12=43\",\"49-43=6\"],\n  [\"79-58=21\",\"86-40=46\"],\n  [\"22-18=4\",\"88-6=82\"],\n  [\"36+23=59\",\"99-73=26\"],\n  [\"3+3=6\",\"72-29=43\"],\n  [\"27+15=42\",\"67-21=46\"],\n  [\"79-72=7\",\"89-84=5\"],\n  [\"57+12=69\",\"39+20=59\"],\n  [\"51+31=82\",\"87-72=15\"],\n  [\"37-37=0\",\"69-65=4\"],\n  [\"24+30=54\",\"46+0=46\"],\n  [\"61-40=21\",\"17+5=22\"],\n  [\"66-52=14\",\"53-35=18\"],\n  [\"1+45=46\",\"0+77=77\"],\n  [\"42-10=32\",\"21+73=94\"],\n  [\"68-35=33\",\"96+2=98\"],\n  [\"11+65=76\",\"32+44=76\"],\n  [\"69+10=79\",\"17+4=21\"],\n  [\"99-44=55\",\"7-3=4\"],\n  [\"11+75=86\",\"28+53=81\"],\n  [\"71-2=69\",\"64+10=74\"],\n  [\"73-34=39\",\"47-5=42\"],\n  [\"16+33=49\",\"12+1=13\"],\n  [\"13+67=80\",\"70-29=41\"],\n  [\"48+3=51\",\"20+46=66\"],\n  [\"87-86=1\",\"14+25=39\"],\n  [\"65-9=56\",\"37+0=37\"],\n  [\"87-41=46\",\"67+22=89\"],\n  [\"95-5=90\",\"59+10=69\"],\n  [\"9+10=19\",\"54-48=6\"],\n  [\"29+68=97\",\"84-20=64\"],\n  [\"95-57=38\",\"85-72=13\"],\n  [\"40+29=69\",\"96-87=9\"],\n  [\"61+36=97\",\"66-49=17\"],\n  [\"81-57=24\",\"42+53=95\"],\n  [\"81-8=73\",\"80-65=15\"],\n  [\"42+10=52\",\"18-2=16\"],\n  [\"63-35=28\",\"52+31=83\"],\n  [\"26+22=48\",\"16+59=75\"],\n  [\"84-53=31\",\"39-22=17\"],\n  [\"17+73=90\",\"67-19=48\"],\n  [\"57-43=14\",\"97-31=66\"],\n  [\"10+12=22\",\"68-34=34\"],\n  [\"39+18=57\",\"6+73=79\"],\n  [\"93-55=38\",\"49-18=31\"],\n  [\"34-8=26\",\"17+37=54\"],\n  [\"80-5=75\",\"56-12=44\"],\n  [\"78-63=15\",\"96-60=36\"],\n  [\"93-74=19\",\"86-6=80\"],\n  [\"7+49=56\",\"9+59=68\"],\n  [\"61+4=65\",\"52+13=65\"],\n  [\"73-3=70\",\"76-1=75\"],\n  [\"10+8=18\",\"69-18=51\"],\n  [\"41+34=75\",\"91-74=17\"],\n  [\"0+36=36\",\"32+58=90\"],\n  [\"33+61=94\",\"57-22=35\"],\n  [\"60+11=71\",\"68+27=95\"],\n  [\"6+62=68\",\"37+57=94\"],\n  [\"3+32=35\",\"80-14=66\"],\n  [\"40-17=23\",\"94-0=94\"],\n  [\"99-76=23\",\"64-56=8\"],\n  [\"29+40=69\",\"50+27=77\"],\n  [\"99-12=87\",\"71-48=23\"],\n  [\"61-21=40\",\"17+70=87\"],\n  [\"73-34=39\",\"73-10=63\"],\n  [\"85-11=74\",\"35-19=16\"],\n  [\"98-4=94\",\"18+30=48\"],\n  [\"16+55=71\",\"72-4=68\"],\n  [\"10+51=61\",\"12+79=91\"],\n  [\"77-70=7\",\"54+23=77\"],\n  [\"26-9=17\",\"54+35=89\"],\n  [\"29-4=25\",\"44+17=61\"],\n  [\"24+7=31\",\"4+82=86\"],\n  [\"36-15=21\",\"21+18=39\"],\n  [\"27-15=12\",\"45-21=24\"],\n  [\"57+39=96\",\"69-26=43\"],\n  [\"60-5=55\",\"39+32=71\"],\n  [\"15+11=26\",\"71+23=94\"],\n  [\"74+13=87\",\"49-49=0\"],\n  [\"66+25=91\",\"58+15=73\"],\n  [\"75-36=39\",\"8+8=16\"],\n  [\"48-20=28\",\"30+55=85\"],\n  [\"72+0=72\",\"26+70=96\"],\n  [\"91-18=73\",\"96-33=63\"]\n];\n\nconst COLS = 5;\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document, found none.\");\n}\n\nconst table = tables.items[0];\n\n// Write each new equation into its cell by (row, col) position. Position-\n// based addressing (rather than text search/replace) is required because a\n// few \"old\" equations are duplicated elsewhere in the table, so matching by\n// text alone cannot tell the repeated cells apart.\nfor (let i = 0; i < REPLACEMENTS.length; i++) {\n  const row = Math.floor(i / COLS);\n  const col = i % COLS;\n  const [, newText] = REPLACEMENTS[i];\n  table.getCell(row, col).value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-answer cells (20 rows x 5 cols) in the single\n# table of the document with new equations, per the diff. Cells are addressed\n# positionally (1-indexed row/col, matching document order) because a handful\n# of \"old\" values repeat, so a naive global Find/Replace would mis-map them.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"56+7=63\"\n$t.Cell(1, 2).Range.Text = \"97-19=78\"\n$t.Cell(1, 3).Range.Text = \"92-20=72\"\n$t.Cell(1, 4).Range.Text = \"18+5=23\"\n$t.Cell(1, 5).Range.Text = \"24+74=98\"\n$t.Cell(2, 1).Range.Text = \"99-34=65\"\n$t.Cell(2, 2).Range.Text = \"89-63=26\"\n$t.Cell(2, 3).Range.Text = \"8+8=16\"\n$t.Cell(2, 4).Range.Text = \"9+3=12\"\n$t.Cell(2, 5).Range.Text = \"74-8=66\"\n$t.Cell(3, 1).Range.Text = \"81-23=58\"\n$t.Cell(3, 2).Range.Text = \"22+0=22\"\n$t.Cell(3, 3).Range.Text = \"43-33=10\"\n$t.Cell(3, 4).Range.Text = \"96-69=27\"\n$t.Cell(3, 5).Range.Text = \"61+23=84\"\n$t.Cell(4, 1).Range.Text = \"32+53=85\"\n$t.Cell(4, 2).Range.Text = \"49-43=6\"\n$t.Cell(4, 3).Range.Text = \"86-40=46\"\n$t.Cell(4, 4).Range.Text = \"88-6=82\"\n$t.Cell(4, 5).Range.Text = \"99-73=26\"\n$t.Cell(5, 1).Range.Text = \"72-29=43\"\n$t.Cell(5, 2).Range.Text = \"67-21=46\"\n$t.Cell(5, 3).Range.Text = \"89-84=5\"\n$t.Cell(5, 4).Range.Text = \"39+20=59\"\n$t.Cell(5, 5).Range.Text = \"87-72=15\"\n$t.Cell(6, 1).Range.Text = \"69-65=4\"\n$t.Cell(6, 2).Range.Text = \"46+0=46\"\n$t.Cell(6, 3).Range.Text = \"17+5=22\"\n$t.Cell(6, 4).Range.Text = \"53-35=18\"\n$t.Cell(6, 5).Range.Text = \"0+77=77\"\n$t.Cell(7, 1).Range.Text = \"21+73=94\"\n$t.Cell(7, 2).Range.Text = \"96+2=98\"\n$t.Cell(7, 3).Range.Text = \"32+44=76\"\n$t.Cell(7, 4).Range.Text = \"17+4=21\"\n$t.Cell(7, 5).Range.Text = \"7-3=4\"\n$t.Cell(8, 1).Range.Text = \"28+53=81\"\n$t.Cell(8, 2).Range.Text = \"64+10=74\"\n$t.Cell(8, 3).Range.Text = \"47-5=42\"\n$t.Cell(8, 4).Range.Text = \"12+1=13\"\n$t.Cell(8, 5).Range.Text = \"70-29=41\"\n$t.Cell(9, 1).Range.Text = \"20+46=66\"\n$t.Cell(9, 2).Range.Text = \"14+25=39\"\n$t.Cell(9, 3).Range.Text = \"37+0=37\"\n$t.Cell(9, 4).Range.Text = \"67+22=89\"\n$t.Cell(9, 5).Range.Text = \"59+10=69\"\n$t.Cell(10, 1).Range.Text = \"54-48=6\"\n$t.Cell(10, 2).Range.Text = \"84-20=64\"\n$t.Cell(10, 3).Range.Text = \"85-72=13\"\n$t.Cell(10, 4).Range.Text = \"96-87=9\"\n$t.Cell(10, 5).Range.Text = \"66-49=17\"\n$t.Cell(11, 1).Range.Text = \"42+53=95\"\n$t.Cell(11, 2).Range.Text = \"80-65=15\"\n$t.Cell(11, 3).Range.Text = \"18-2=16\"\n$t.Cell(11, 4).Range.Text = \"52+31=83\"\n$t.Cell(11, 5).Range.Text = \"16+59=75\"\n$t.Cell(12, 1).Range.Text = \"39-22=17\"\n$t.Cell(12, 2).Range.Text = \"67-19=48\"\n$t.Cell(12, 3).Range.Text = \"97-31=66\"\n$t.Cell(12, 4).Range.Text = \"68-34=34\"\n$t.Cell(12, 5).Range.Text = \"6+73=79\"\n$t.Cell(13, 1).Range.Text = \"49-18=31\"\n$t.Cell(13, 2).Range.Text = \"17+37=54\"\n$t.Cell(13, 3).Range.Text = \"56-12=44\"\n$t.Cell(13, 4).Range.Text = \"96-60=36\"\n$t.Cell(13, 5).Range.Text = \"86-6=80\"\n$t.Cell(14, 1).Range.Text = \"9+59=68\"\n$t.Cell(14, 2).Range.Text = \"52+13=65\"\n$t.Cell(14, 3).Range.Text = \"76-1=75\"\n$t.Cell(14, 4).Range.Text = \"69-18=51\"\n$t.Cell(14, 5).Range.Text = \"91-74=17\"\n$t.Cell(15, 1).Range.Text = \"32+58=90\"\n$t.Cell(15, 2).Range.Text = \"57-22=35\"\n$t.Cell(15, 3).Range.Text = \"68+27=95\"\n$t.Cell(15, 4).Range.Text = \"37+57=94\"\n$t.Cell(15, 5).Range.Text = \"80-14=66\"\n$t.Cell(16, 1).Range.Text = \"94-0=94\"\n$t.Cell(16, 2).Range.Text = \"64-56=8\"\n$t.Cell(16, 3).Range.Text = \"50+27=77\"\n$t.Cell(16, 4).Range.Text = \"71-48=23\"\n$t.Cell(16, 5).Range.Text = \"17+70=87\"\n$t.Cell(17, 1).Range.Text = \"73-10=63\"\n$t.Cell(17, 2).Range.Text = \"35-19=16\"\n$t.Cell(17, 3).Range.Text = \"18+30=48\"\n$t.Cell(17, 4).Range.Text = \"72-4=68\"\n$t.Cell(17, 5).Range.Text = \"12+79=91\"\n$t.Cell(18, 1).Range.Text = \"54+23=77\"\n$t.Cell(18, 2).Range.Text = \"54+35=89\"\n$t.Cell(18, 3).Range.Text = \"44+17=61\"\n$t.Cell(18, 4).Range.Text = \"4+82=86\"\n$t.Cell(18, 5).Range.Text = \"21+18=39\"\n$t.Cell(19, 1).Range.Text = \"45-21=24\"\n$t.Cell(19, 2).Range.Text = \"69-26=43\"\n$t.Cell(19, 3).Range.Text = \"39+32=71\"\n$t.Cell(19, 4).Range.Text = \"71+23=94\"\n$t.Cell(19, 5).Range.Text = \"49-49=0\"\n$t.Cell(20, 1).Range.Text = \"58+15=73\"\n$t.Cell(20, 2).Range.Text = \"8+8=16\"\n$t.Cell(20, 3).Range.Text = \"30+55=85\"\n$t.Cell(20, 4).Range.Text = \"26+70=96\"\n$t.Cell(20, 5).Range.Text = \"96-33=63\"\n"}
